$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Date/time number format used by column A (matches existing style s="2")
$dateFmt = "YYYY-MM-DD HH:MM:SS"

# Helper: write a truly-empty text cell (matches an empty inline string <is/>
# in the source data) without leaving a stray quote-prefix style behind.
function Set-BlankText($row, $col) {
    $c = $ws.Cells.Item($row, $col)
    $c.Value = "'"
    $c.Style = "Normal"
}

# --- Row 2 ---
$ws.Cells.Item(2,1).Value = 45862.45833333334
Set-BlankText 2 2
Set-BlankText 2 3
$ws.Cells.Item(2,4).Value = 16.44
$ws.Cells.Item(2,5).Value = 83.03
$ws.Cells.Item(2,6).Value = 648.1900000000001
$ws.Cells.Item(2,7).Value = 13.25
$ws.Cells.Item(2,8).Value = "SE"
$ws.Cells.Item(2,9).Value = 0
Set-BlankText 2 10

# --- Row 3 ---
$ws.Cells.Item(3,1).Value = 45862.5
Set-BlankText 3 2
Set-BlankText 3 3
$ws.Cells.Item(3,4).Value = 17.65
$ws.Cells.Item(3,5).Value = 79.31999999999999
$ws.Cells.Item(3,6).Value = 613.04
$ws.Cells.Item(3,7).Value = 13.12
$ws.Cells.Item(3,8).Value = "ESE"
$ws.Cells.Item(3,9).Value = 0
Set-BlankText 3 10

# --- Row 4 ---
$ws.Cells.Item(4,1).Value = 45862.54166666666
Set-BlankText 4 2
Set-BlankText 4 3
$ws.Cells.Item(4,4).Value = 18.51
$ws.Cells.Item(4,5).Value = 76.68000000000001
$ws.Cells.Item(4,6).Value = 327.83
$ws.Cells.Item(4,7).Value = 14.07
$ws.Cells.Item(4,8).Value = "ESE"
$ws.Cells.Item(4,9).Value = 0
Set-BlankText 4 10

# --- Row 5 ---
$ws.Cells.Item(5,1).Value = 45862.58333333334
Set-BlankText 5 2
Set-BlankText 5 3
$ws.Cells.Item(5,4).Value = 19.33
$ws.Cells.Item(5,5).Value = 74.44
$ws.Cells.Item(5,6).Value = 87.89
$ws.Cells.Item(5,7).Value = 13.54
$ws.Cells.Item(5,8).Value = "ESE"
$ws.Cells.Item(5,9).Value = 0
Set-BlankText 5 10

# --- Row 6 (new, mostly blank) ---
$ws.Cells.Item(6,1).Value = 45862.625
$ws.Cells.Item(6,1).NumberFormat = $dateFmt
Set-BlankText 6 2
Set-BlankText 6 3
Set-BlankText 6 4
Set-BlankText 6 5
Set-BlankText 6 6
Set-BlankText 6 7
Set-BlankText 6 8
Set-BlankText 6 9
Set-BlankText 6 10

# --- Row 7 (new, mostly blank) ---
$ws.Cells.Item(7,1).Value = 45862.66666666666
$ws.Cells.Item(7,1).NumberFormat = $dateFmt
Set-BlankText 7 2
Set-BlankText 7 3
Set-BlankText 7 4
Set-BlankText 7 5
Set-BlankText 7 6
Set-BlankText 7 7
Set-BlankText 7 8
Set-BlankText 7 9
Set-BlankText 7 10

# --- Row 8 (new, full data - latest reading) ---
$ws.Cells.Item(8,1).Value = 45862.58333333334
$ws.Cells.Item(8,1).NumberFormat = $dateFmt
$ws.Cells.Item(8,2).Value = 2025
$ws.Cells.Item(8,3).Value = 30
$ws.Cells.Item(8,4).Value = 19.33
$ws.Cells.Item(8,5).Value = 74.44
$ws.Cells.Item(8,6).Value = 87.89
$ws.Cells.Item(8,7).Value = 13.54
$ws.Cells.Item(8,8).Value = "ESE"
$ws.Cells.Item(8,9).Value = 0
$ws.Cells.Item(8,10).Value = "16:38:57"
